$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values for row 2 and row 3
$ws.Range("B2").Value = 349
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 120

# Remove rows 4 and 5 entirely (last cases trimmed from the control points table)
$ws.Range("A4:B5").Delete()
